$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Add the new "Cohort" worksheet as the last tab in the workbook
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Cohort"

# ------------------------------------------------------------------
# Column widths
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 118.21875
$ws.Columns.Item(2).ColumnWidth = 35.6640625

# ------------------------------------------------------------------
# Formatting: column A data rows get Arial/black font, left/center
# aligned, wrapped text (re-use the existing Arial style already in
# the workbook, then tweak indent + wrap).
# ------------------------------------------------------------------
$fmtSource = $wb.Worksheets.Item("Freshman Profile_Percentile").Range("A2")
$questionRange = $ws.Range("A2:A19")
$fmtSource.Copy()
$questionRange.PasteSpecial(-4122)  # xlPasteFormats
$questionRange.IndentLevel = 0
$questionRange.WrapText = $true
$excel.CutCopyMode = $false

# Header cell (A1) just wraps text
$ws.Range("A1").WrapText = $true

# Column B (answers) stored as text
$ws.Range("B1:B19").NumberFormat = "@"

# ------------------------------------------------------------------
# Header row
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Question"
$ws.Range("B1").Value = "Answer"

# ------------------------------------------------------------------
# Data rows
# ------------------------------------------------------------------
$questions = @(
    "How many students in the initial cohort were recipients of a Federal Pell Grant?",
    "How many of those students in the initial cohort were recipients of a Subsidized Stafford Loan but not a Pell Grant?",
    "How many of those students in the initial cohort did not receive either a Pell Grant or a subsidized Stafford Loan?",
    "How many exempted student were recipients of a Federal Pell Grant?",
    "How many of exempted students were recipients of a Subsidized Stafford Loan but not a Pell Grant?",
    "How many of exempted students did not receive a Pell Grant or a subsidized Stafford Loan?",
    "How many students in the final cohort were recipients of a Federal Pell Grant?",
    "How many students in the final cohort were recipients of a Subsidized Stafford Loan but not a Pell Grant?",
    "How many students in the final cohort did not receive a Pell Grant or a subsidized Stafford Loan?",
    "How many students in the initial cohort completed the program within four years and were recipients of a Federal Pell Grant?",
    "How many students in the initial cohort who completed the program within four years were recipients of a Subsidized Stafford Loan but not a Pell Grant?",
    "How many students in the initial cohort who completed the program within four years but did not receive a Pell Grant or a subsidized Stafford Loan?",
    "How many students  in the initial cohort who completed the program in between four and five years were recipients of a Federal Pell Grant?",
    "How many students  in the initial cohort who completed the program in between four and five years were recipients of a Subsidized Stafford Loan but not a Pell Grant?",
    "How many students in the initial cohort who completed the program in between four and five years did not receive a Pell Grant or a subsidized Stafford Loan?",
    "How many students who completed the program between five and six years were recipients of a Federal Pell Grant?",
    "How many students who completed the program between five and six years were recipients of a Subsidized Stafford Loan but not a Pell Grant?",
    "How many students who completed the program between five and six years did not receive a Pell Grant or a subsidized Stafford Loan?"
)

$answers = @("83", "199", "300", "0", "0", "1", "83", "199", "299", "61", "144", "220", "7", "9", "33", "1", "5", "3")

for ($i = 0; $i -lt $questions.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $questions[$i]
    $ws.Cells.Item($r, 2).Value = $answers[$i]
}

# Rows whose question text wraps onto two lines in the source workbook
$tallRows = @(12, 13, 14, 15, 16, 18, 19)
foreach ($r in $tallRows) {
    $ws.Rows.Item($r).RowHeight = 27.6
}

# ------------------------------------------------------------------
# Page setup
# ------------------------------------------------------------------
$ws.PageSetup.Orientation = 1  # xlPortrait

# ------------------------------------------------------------------
# View state: Cohort becomes the active/selected sheet & tab
# ------------------------------------------------------------------
$freshman = $wb.Worksheets.Item("Freshman Profile_Percentile")
$freshman.Range("A43").Select()

$ws.Activate()
$ws.Range("E15").Select()
